# Update the marksheet: correct the "Corr/total marks" figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: value in B11 changes from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row: value in B12 changes from 78 to 130
$ws.Range("B12").Value = 130

# "Total" row: correct/total marks text in E12 changes from "77/84" to "130/140"
$ws.Range("E12").Value = "130/140"
